$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above row 51 (currently "DUALCAR"), shifting DUALCAR and
# everything below it down by one row. DXLINE sorts alphabetically between
# DUALCAR and ERASENOSE.
$row = $ws.Rows.Item(51)
$row.Insert()

$ws.Cells.Item(51, 1).Value = "DXLINE"
$ws.Cells.Item(51, 2).Value = "Draws an XLINE at a given or selected station"
